$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-matrix cell values resulting from additional
# simulated games (more games run -> updated empirical frequencies).
    $ws.Range("B2").Value = 0.1965811965811966
    $ws.Range("C2").Value = 0.5085470085470085
    $ws.Range("J2").Value = 0.02136752136752137
    $ws.Range("P2").Value = 0.141025641025641
    $ws.Range("S2").Value = 0.1324786324786325
    $ws.Range("B3").Value = 0.008264462809917356
    $ws.Range("C3").Value = 0.01652892561983471
    $ws.Range("J3").Value = 0.04132231404958678
    $ws.Range("P3").Value = 0.7603305785123967
    $ws.Range("S3").Value = 0.1735537190082645
    $ws.Range("J4").Value = 0.05
    $ws.Range("P4").Value = 0.55
    $ws.Range("S4").Value = 0.4
    $ws.Range("B6").Value = 0.08
    $ws.Range("D6").Value = 0.01714285714285714
    $ws.Range("F6").Value = 0.04
    $ws.Range("J6").Value = 0.1942857142857143
    $ws.Range("O6").Value = 0.01142857142857143
    $ws.Range("Q6").Value = 0.2057142857142857
    $ws.Range("R6").Value = 0.05142857142857143
    $ws.Range("S6").Value = 0.4
    $ws.Range("B7").Value = 0.0650887573964497
    $ws.Range("D7").Value = 0.005917159763313609
    $ws.Range("F7").Value = 0.04142011834319527
    $ws.Range("J7").Value = 0.136094674556213
    $ws.Range("O7").Value = 0.02958579881656805
    $ws.Range("Q7").Value = 0.2130177514792899
    $ws.Range("R7").Value = 0.1005917159763314
    $ws.Range("S7").Value = 0.408284023668639
    $ws.Range("B8").Value = 0.0594900849858357
    $ws.Range("D8").Value = 0.0226628895184136
    $ws.Range("F8").Value = 0.0594900849858357
    $ws.Range("J8").Value = 0.08781869688385269
    $ws.Range("O8").Value = 0.0339943342776204
    $ws.Range("Q8").Value = 0.1784702549575071
    $ws.Range("R8").Value = 0.08781869688385269
    $ws.Range("S8").Value = 0.4702549575070821
    $ws.Range("B9").Value = 0.08333333333333333
    $ws.Range("D9").Value = 0.01785714285714286
    $ws.Range("F9").Value = 0.07738095238095238
    $ws.Range("J9").Value = 0.1071428571428571
    $ws.Range("O9").Value = 0.03571428571428571
    $ws.Range("Q9").Value = 0.1488095238095238
    $ws.Range("R9").Value = 0.09523809523809523
    $ws.Range("S9").Value = 0.4345238095238095
    $ws.Range("B10").Value = 0.1087751371115174
    $ws.Range("D10").Value = 0.02285191956124314
    $ws.Range("F10").Value = 0.06764168190127971
    $ws.Range("J10").Value = 0.1124314442413163
    $ws.Range("O10").Value = 0.01188299817184644
    $ws.Range("Q10").Value = 0.2138939670932358
    $ws.Range("R10").Value = 0.06398537477148081
    $ws.Range("S10").Value = 0.3985374771480805
    $ws.Range("G11").Value = 0.1680327868852459
    $ws.Range("J11").Value = 0.07377049180327869
    $ws.Range("K11").Value = 0.1967213114754098
    $ws.Range("L11").Value = 0.5491803278688525
    $ws.Range("S11").Value = 0.01229508196721311
    $ws.Range("G12").Value = 0.7372262773722628
    $ws.Range("J12").Value = 0.1678832116788321
    $ws.Range("K12").Value = 0.0072992700729927
    $ws.Range("L12").Value = 0.05109489051094891
    $ws.Range("S12").Value = 0.0364963503649635
    $ws.Range("G13").Value = 0.6578947368421053
    $ws.Range("J13").Value = 0.2368421052631579
    $ws.Range("S13").Value = 0.1052631578947368
    $ws.Range("F15").Value = 0.01
    $ws.Range("H15").Value = 0.12
    $ws.Range("I15").Value = 0.095
    $ws.Range("J15").Value = 0.36
    $ws.Range("K15").Value = 0.06
    $ws.Range("M15").Value = 0.01
    $ws.Range("O15").Value = 0.095
    $ws.Range("S15").Value = 0.25
    $ws.Range("F16").Value = 0.01438848920863309
    $ws.Range("H16").Value = 0.1294964028776978
    $ws.Range("I16").Value = 0.07194244604316546
    $ws.Range("J16").Value = 0.5035971223021583
    $ws.Range("K16").Value = 0.1151079136690648
    $ws.Range("M16").Value = 0.01438848920863309
    $ws.Range("O16").Value = 0.05755395683453238
    $ws.Range("S16").Value = 0.09352517985611511
    $ws.Range("F17").Value = 0.01503759398496241
    $ws.Range("H17").Value = 0.1528822055137845
    $ws.Range("I17").Value = 0.05764411027568922
    $ws.Range("J17").Value = 0.4385964912280702
    $ws.Range("K17").Value = 0.09022556390977443
    $ws.Range("M17").Value = 0.02506265664160401
    $ws.Range("O17").Value = 0.06516290726817042
    $ws.Range("S17").Value = 0.1553884711779449
    $ws.Range("F18").Value = 0.01398601398601399
    $ws.Range("H18").Value = 0.1468531468531468
    $ws.Range("I18").Value = 0.0979020979020979
    $ws.Range("J18").Value = 0.3846153846153846
    $ws.Range("K18").Value = 0.07692307692307693
    $ws.Range("M18").Value = 0.02797202797202797
    $ws.Range("O18").Value = 0.08391608391608392
    $ws.Range("S18").Value = 0.1678321678321678
    $ws.Range("F19").Value = 0.0176619007569386
    $ws.Range("H19").Value = 0.1976450798990748
    $ws.Range("I19").Value = 0.08662741799831791
    $ws.Range("J19").Value = 0.3734230445752734
    $ws.Range("K19").Value = 0.1000841042893188
    $ws.Range("M19").Value = 0.0159798149705635
    $ws.Range("N19").Value = 0.004205214465937763
    $ws.Range("O19").Value = 0.05971404541631623
    $ws.Range("S19").Value = 0.144659377628259
